$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hospital-name list was refreshed: "VS Hospitals" was dropped, the
# remaining entries shifted up, and a new hospital was appended at the end.
# Row 10 ("VS Hospitals") now reads the old row-11 text, and row 11 gets the
# newly added hospital.
$ws.Range("A10").Value = "Jayam Hospital & GFC Fertility"
$ws.Range("A11").Value = "Maaya Speciality Hospitals"

# A second column was introduced next to the hospital names (multiple
# browsing / per-hospital link feature), which is why the sheet's used
# range now extends from A1:A11 to A1:B11.
$ws.Range("B1").Value = "Hospital Link"
